# Update cryptos list: refresh Price/Volume(1h) values and fix a
# Stacks/EnergySwap row-order swap, per the Thu Feb 15 20:51:03 UTC 2024
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "51.808.06"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.08%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "2.835.45"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.57%  "
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "352.63"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "113.65"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.55%  "
$ws.Cells.Item(7, 5).Value = "  +3.81%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.05%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.600"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +4.03%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "41.63"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.58%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0851"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.98%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "19.96"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -1.20%  "
$ws.Cells.Item(13, 5).Value = "  +1.45%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "7.72"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.78%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "3.280.53"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.67%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "2.824.36"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.07%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.896"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.57%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "51.723.58"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.01%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +7.32%  "
$ws.Cells.Item(20, 5).Value = "  -1.85%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "13.50"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.31%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0994"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.09%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "270.70"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.83%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "69.73"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.11%  "
$ws.Cells.Item(25, 5).Value = "  +3.59%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "26.71"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.40%  "
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "10.29"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.94%  "
$ws.Cells.Item(29, 5).Value = "  +1.13%  "
$ws.Cells.Item(30, 5).Value = "  -1.43%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "50.67"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.12%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "33.90"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -3.43%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.0449"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +27.66%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "5.81"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.27%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.0827"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.17%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.06%  "
$ws.Cells.Item(37, 5).Value = "  -0.33%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "3.22"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.59%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "4.88"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.61%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "18.00"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -5.23%  "
$ws.Cells.Item(41, 2).Value = "EnergySwap"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "23.69"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.02%  "
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +4.65%  "
$ws.Cells.Item(43, 5).Value = "  +1.04%  "
$ws.Cells.Item(44, 5).Value = "  -1.25%  "
$ws.Cells.Item(45, 5).Value = "  +0.22%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "2.079.90"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.50%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "3.33"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.24%  "
$ws.Cells.Item(48, 5).Value = "  +3.61%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "5.70"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.93%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.935"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +6.78%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "60.80"
$c.Style = "Normal"
